$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of portfolio data (row 39) after the last existing row (38).
# Column A holds a date-like string that must stay plain text (matching the
# other rows), so format it as text before assigning, then reset the style
# back to "Normal" so no extra/unused style gets attached to the cell.
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "2025-09-23"
$ws.Range("A39").Style = "Normal"

$ws.Range("B39").Value = 58.75
$ws.Range("C39").Value = 701.3499755859375
$ws.Range("D39").Value = 338.3500061035156
